# CodeSystem-MedioDeLlegada: publish new version (status -> active, bump date,
# and record Case Sensitive = true)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: bump published timestamp
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true
# Use a leading apostrophe so Excel stores this as text (matching the
# column's existing string type) rather than auto-converting to a boolean,
# then restore the original (non "quote prefix") cell formatting by pasting
# just the format from an untouched neighboring cell.
$ws.Range("B17").Formula = "'true"
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
